$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# Fixed comparison counter values for Knuth Morris Pratt (column F), rows 16-23
$ws.Range("F16").Value = 18
$ws.Range("F17").Value = 190
$ws.Range("F18").Value = 1900
$ws.Range("F19").Value = 9000
$ws.Range("F20").Value = 633
$ws.Range("F21").Value = 695
$ws.Range("F22").Value = 2070
$ws.Range("F23").Value = 9614
